$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.846373975276947
$ws.Range("B1").Value = 1.407367467880249
$ws.Range("C1").Value = 3.790290355682373
$ws.Range("D1").Value = 2.688522577285767
$ws.Range("E1").Value = 1.617789387702942
